$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# Correct the "FREESTATE" typo to "FREE STATE" everywhere it is used
$ws2.Range("B7").Value = "FREE STATE"
$ws2.Range("B12").Value = "FREE STATE"
$ws2.Range("B17").Value = "FREE STATE"

# Add the missing FREE STATE / EASTERN CAPE data row
$ws2.Range("A18").Value = 16
$ws2.Range("B18").Value = "FREE STATE"
$ws2.Range("C18").Value = "EASTERN CAPE"
$ws2.Range("D18").Value = 0
$ws2.Range("E18").Value = 0
$ws2.Range("F18").Value = 0

# Switch the active sheet/selection to the 9-supply sheet
$ws2.Activate()
$ws2.Range("D24").Select()
